$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# -----------------------------------------------------------------
# 1) Header cells: "Pass/Fail" -> "Result" for each test-block header
#    (rows 14, 18, 73, 80, 87, 94, column B)
# -----------------------------------------------------------------
$headerRows = @(14, 18, 73, 80, 87, 94)
foreach ($hr in $headerRows) {
    $ws.Range("B$hr").Value = "Result"
}

# -----------------------------------------------------------------
# 2) Simple "Done" rows 15 and 19: clear B, put formula in C
# -----------------------------------------------------------------
$simpleRows = @(15, 19)
foreach ($row in $simpleRows) {
    $ws.Range("B$row").ClearContents()
    $ws.Range("C$row").Formula = '=IF(ISBLANK(B' + $row + '),"Not Done",B' + $row + ')'
}

# -----------------------------------------------------------------
# 3) Grouped "Done" rows (74-77, 81-84, 88-91, 95-98):
#    clear B, remove the now-unused J marker cell, put formula in C
#    (first row of each group gets its own formula, the remaining
#    three share one formula across the C<first+1>:C<last> range so
#    Excel emits a shared formula, matching the source workbook)
# -----------------------------------------------------------------
$groupStarts = @(74, 81, 88, 95)
foreach ($start in $groupStarts) {
    $last = $start + 3

    for ($r = $start; $r -le $last; $r++) {
        $ws.Range("B$r").ClearContents()
        $ws.Range("J$r").ClearContents()
    }

    $ws.Range("C$start").Formula = '=IF(ISBLANK(B' + $start + '),"Not Done",B' + $start + ')'

    $secondRow = $start + 1
    $ws.Range("C${secondRow}:C${last}").Formula = '=IF(ISBLANK(B' + $secondRow + '),"Not Done",B' + $secondRow + ')'

    # ---------------------------------------------------------
    # Conditional formatting: the four single-cell "highlight"
    # rules (ISBLANK / Fail / Pass) that used to live on each
    # row's C cell are consolidated onto the whole C<start>:C<last>
    # range.
    # ---------------------------------------------------------
    for ($r = $start; $r -le $last; $r++) {
        $ws.Range("C$r").FormatConditions.Delete()
    }

    $rng = $ws.Range("C$start`:C$last")

    $fcBlank = $rng.FormatConditions.Add(9, 0, "highlight")
    $fcBlank.Formula1 = "=ISBLANK(C$start)"
    $fcBlank.Text = "highlight"
    $fcBlank.Interior.Color = 10086143

    $fcFail = $rng.FormatConditions.Add(9, 0, "highlight")
    $fcFail.Formula1 = '=C' + $start + '="Fail"'
    $fcFail.Text = "highlight"
    $fcFail.Font.Color = 255

    $fcPass = $rng.FormatConditions.Add(9, 0, "highlight")
    $fcPass.Formula1 = '=C' + $start + '="Pass"'
    $fcPass.Text = "highlight"
    $fcPass.Font.Color = 5287936
}
